$d = $word.ActiveDocument
$d.Content.Find.Execute("loop.index + table_page_length - 1 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 ", 2)
